$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "H 72" record (row 2). This shifts all subsequent rows up by
# one, so the last row (previously row 63, "SC 232") moves to row 62 and
# the used range shrinks from A1:F63 to A1:F62.
$ws.Rows.Item(2).Delete()
